# Cobalt Users workbook update
# - Rename Sheet2 -> Emails and populate it with ANZ tester creds
# - Add a block of new QA users (Search/FFH/FrontEnd/Url/Linking/Login/Cpet) to Users sheet,
#   each with a hard-coded warning note, Locked="N" and a mailto hyperlink on the Email column
# - Trim Sheet3 back down to a single "Y" flag row
# - Restore the previous selection / column-width look & feel

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Users sheet - append the new test accounts
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Users")

$usernames = @(
    "SearchOpenWebUser1",
    "FFHUser1", "FFHUser2", "FFHUser3", "FFHUser4",
    "FrontEndUser1", "FrontEndUser2", "FrontEndUser3", "FrontEndUser4", "FrontEndUser5",
    "FrontEndUser6", "FrontEndUser7", "FrontEndUser8", "FrontEndUser9", "FrontEndUser10",
    "UrlUser1", "UrlUser2", "UrlUser3",
    "LinkingUser1",
    "LoginUser1", "LoginUser2", "LoginUser3", "LoginUser4", "LoginUser5", "LoginUser6", "LoginUser7",
    "CpetUser1", "CpetUser2"
)

$emails = @(
    "SearchOpenWeb@mailinator.com ",
    "FFHUser1@mailinator.com ", "FFHUser2@mailinator.com", "FFHUser3@mailinator.com", "FFHUser4@mailinator.com",
    "FrontEndUser1@mailinator.com", "FrontEndUser2@mailinator.com", "FrontEndUser3@mailinator.com", "FrontEndUser4@mailinator.com", "FrontEndUser5@mailinator.com",
    "FrontEndUser6@mailinator.com", "FrontEndUser7@mailinator.com", "FrontEndUser8@mailinator.com", "FrontEndUser9@mailinator.com", "FrontEndUser10@mailinator.com",
    "UrlUser1@mailinator.com", "UrlUser2@mailinator.com", "UrlUser3@mailinator.com",
    "LinkingUser1@mailinator.com ",
    "LoginUser1@mailinator.com ", "LoginUser2@mailinator.com ", "LoginUser3@mailinator.com ", "LoginUser4@mailinator.com ", "LoginUser5@mailinator.com ", "LoginUser6@mailinator.com ", "LoginUser7@mailinator.com ",
    "CpetUser1@mailinator.com ", "CpetUser2@mailinator.com "
)

$firstRow = 53
$note = "THIS IS IN USE 24/7 - DO NOT USE!"

# Column A - all usernames first (matches shared-string build order of the source file)
for ($i = 0; $i -lt $usernames.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $usernames[$i]
}

# Column B - reuses the existing "Password1" shared string
for ($i = 0; $i -lt $usernames.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = "Password1"
}

# Column G - all emails next
for ($i = 0; $i -lt $usernames.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 7).Value = $emails[$i]
}

# hyperlinks were recorded in this (non-sequential) order in the source file
$hyperlinkRowOrder = @(55,54,56,57,58,59,60,61,62,63,64,65,66,67,68,69,53,70,71,72,73,74,75,76,77,78,79,80)
foreach ($r in $hyperlinkRowOrder) {
    $idx = $r - $firstRow
    $addr = "mailto:" + $emails[$idx]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 7), $addr) | Out-Null
}

# Column E (the warning note) and F ("N") last
for ($i = 0; $i -lt $usernames.Length; $i++) {
    $r = $firstRow + $i

    $ws.Cells.Item($r, 5).Value = $note
    $ws.Cells.Item($r, 6).Value = "N"

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(10).LineStyle = 1

    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Borders.Item(7).LineStyle = 1
    $fCell.Borders.Item(10).LineStyle = 1
}

# four trailing blank (but bordered) rows under the new block
for ($r = 81; $r -le 84; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(10).LineStyle = 1
}

$ws.Columns("A").AutoFit() | Out-Null
$ws.Columns("E").AutoFit() | Out-Null

$ws.Activate()
$ws.Range("C80").Select()

# ---------------------------------------------------------------------------
# 2. Sheet2 -> Emails
# ---------------------------------------------------------------------------
$wsEmails = $wb.Worksheets.Item("Sheet2")
$wsEmails.Name = "Emails"

$wsEmails.Cells.Item(1, 1).Value = "Email"
$wsEmails.Cells.Item(1, 2).Value = "Password"
$wsEmails.Cells.Item(2, 1).Value = "tr-anz-tester1@yandex.com"
$wsEmails.Cells.Item(2, 2).Value = "tranztest"
$wsEmails.Cells.Item(3, 1).Value = "tr-anz-tester2@yandex.com"
$wsEmails.Cells.Item(3, 2).Value = "tranztest"

$wsEmails.Columns("A").ColumnWidth = 26.28515625
$wsEmails.Columns("B").ColumnWidth = 13.85546875

$wsEmails.Range("A1:B3").Select()

# ---------------------------------------------------------------------------
# 3. Sheet3 - drop the duplicated rows, keep a single "Y"
# ---------------------------------------------------------------------------
$wsFlags = $wb.Worksheets.Item("Sheet3")
$wsFlags.Rows("3:6").Delete()

# ---------------------------------------------------------------------------
# Re-activate Users so it stays the selected tab on save
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C80").Select()
